$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 160.81818
$ws.Range("I5").Value = 160.81818
$ws.Range("K5").Value = 160.81818
$ws.Range("M5").Value = -45.81818000000001
$ws.Range("H9").Value = 158.45454
$ws.Range("I9").Value = 206.28572
$ws.Range("J9").Value = 74.75
$ws.Range("K9").Value = 206.28572
$ws.Range("L9").Value = 74.75
$ws.Range("M9").Value = -37.28572
$ws.Range("N9").Value = -412.75
$ws.Range("H19").Value = 449.83334
$ws.Range("J19").Value = 524.75
$ws.Range("L19").Value = 524.75
$ws.Range("N19").Value = -874.75
$ws.Range("H33").Value = 1092.3
$ws.Range("I33").Value = 115.125
$ws.Range("K33").Value = 115.125
$ws.Range("M33").Value = 113.875
$ws.Range("H43").Value = 5499.5
$ws.Range("I43").Value = 6999.5
$ws.Range("J43").Value = 3999.5
$ws.Range("K43").Value = 6999.5
$ws.Range("L43").Value = 3999.5
$ws.Range("M43").Value = -6930.5
$ws.Range("N43").Value = -4137.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 2016
$ws.Range("I50").Value = 2524
$ws.Range("J50").Value = 1000
$ws.Range("K50").Value = 2524
$ws.Range("L50").Value = 1000
$ws.Range("M50").Value = -1810
$ws.Range("N50").Value = -2428
$ws.Range("H101").Value = 56867
$ws.Range("J101").Value = 56867
$ws.Range("L101").Value = 56867
$ws.Range("N101").Value = -63357
$ws.Range("H132").Value = 1425.1428
$ws.Range("I132").Value = 1425.1428
$ws.Range("K132").Value = 4275.428400000001
$ws.Range("M132").Value = -1745.428400000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3041.8333
$ws.Range("I107").Value = 1687.75
$ws.Range("K107").Value = 1687.75
$ws.Range("M107").Value = 232.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 311.26086
$ws.Range("I7").Value = 432.72726
$ws.Range("J7").Value = 199.91667
$ws.Range("K7").Value = 432.72726
$ws.Range("L7").Value = 199.91667
$ws.Range("M7").Value = -319.72726
$ws.Range("N7").Value = -425.91667
$ws.Range("H19").Value = 486.1111
$ws.Range("I19").Value = 534.375
$ws.Range("K19").Value = 534.375
$ws.Range("M19").Value = -364.375
$ws.Range("H24").Value = 486.1111
$ws.Range("I24").Value = 534.375
$ws.Range("K24").Value = 534.375
$ws.Range("M24").Value = -364.375
$ws.Range("H107").Value = 1242.9
$ws.Range("I107").Value = 775.7143
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 775.7143
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = 1144.2857
$ws.Range("N107").Value = -6173

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.541668
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 55.615383
$ws.Range("K2").Value = 228
$ws.Range("L2").Value = 333.692298
$ws.Range("M2").Value = -115
$ws.Range("N2").Value = -559.6922979999999
$ws.Range("H11").Value = 3333732.5
$ws.Range("I11").Value = 5000477
$ws.Range("J11").Value = 243.4
$ws.Range("K11").Value = 15001431
$ws.Range("L11").Value = 730.2
$ws.Range("M11").Value = -15001291
$ws.Range("N11").Value = -1010.2
$ws.Range("H80").Value = 4968.7
$ws.Range("I80").Value = 2562.3333
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 7686.999899999999
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -6750.999899999999
$ws.Range("N80").Value = -19872
$ws.Range("H83").Value = 4968.7
$ws.Range("I83").Value = 2562.3333
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 23060.9997
$ws.Range("L83").Value = 54000
$ws.Range("M83").Value = -18380.9997
$ws.Range("N83").Value = -63360
$ws.Range("H97").Value = 780.8333
$ws.Range("I97").Value = 728.6667
$ws.Range("J97").Value = 833
$ws.Range("K97").Value = 2186.0001
$ws.Range("L97").Value = 2499
$ws.Range("M97").Value = -1690.0001
$ws.Range("N97").Value = -3491
$ws.Range("H132").Value = 1736.25
$ws.Range("I132").Value = 975
$ws.Range("J132").Value = 2497.5
$ws.Range("K132").Value = 8775
$ws.Range("L132").Value = 22477.5
$ws.Range("M132").Value = -6245
$ws.Range("N132").Value = -27537.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9193
$ws.Range("I46").Value = 9193
$ws.Range("K46").Value = 9193
$ws.Range("M46").Value = -9037
$ws.Range("H80").Value = 3800
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 3800
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = $null
$ws.Range("H102").Value = 1798.7222
$ws.Range("I102").Value = 1859.2941
$ws.Range("K102").Value = 1859.2941
$ws.Range("M102").Value = -237.2941000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2559.818
$ws.Range("J22").Value = 2876
$ws.Range("L22").Value = 2876
$ws.Range("N22").Value = -3466
$ws.Range("H27").Value = 2559.818
$ws.Range("J27").Value = 2876
$ws.Range("L27").Value = 2876
$ws.Range("N27").Value = -3090
$ws.Range("H40").Value = 9225
$ws.Range("I40").Value = 9906.3125
$ws.Range("J40").Value = 6499.75
$ws.Range("K40").Value = 9906.3125
$ws.Range("L40").Value = 6499.75
$ws.Range("M40").Value = -9770.3125
$ws.Range("N40").Value = -6771.75
$ws.Range("H46").Value = 4071.2354
$ws.Range("J46").Value = 4824
$ws.Range("L46").Value = 4824
$ws.Range("N46").Value = -5200
$ws.Range("H56").Value = 18525.5
$ws.Range("I56").Value = 18525.5
$ws.Range("K56").Value = 18525.5
$ws.Range("M56").Value = -17834.5
$ws.Range("H132").Value = 10668
$ws.Range("I132").Value = 10668
$ws.Range("K132").Value = 32004
$ws.Range("M132").Value = -29474

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 5156.6
$ws.Range("I55").Value = 5156.6
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 5156.6
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -4879.6
$ws.Range("N55").Value = $null
$ws.Range("H132").Value = 2277.3845
$ws.Range("I132").Value = 2203.7778
$ws.Range("K132").Value = 6611.3334
$ws.Range("M132").Value = -4081.3334
